$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# toi_source's class changes from "character" to "categorical"
$ws.Range("C6").Value = "categorical"

# Leave the active selection where the user ended up after editing
$null = $ws.Range("B17").Select()
